$wb = $excel.ActiveWorkbook

# "Multi-storm shps" sheet: clear the leftover empty/styled cells in row 2 (A2:B2),
# keeping D2's note text intact.
$wsMulti = $wb.Worksheets.Item("Multi-storm shps")
$wsMulti.Range("A2:B2").Clear()
$wsMulti.Range("A2:C6").Select()

# "Single-storm shps" sheet: clear the leftover empty/styled cells and the stray
# storm_rank value of 1 in row 2 (A2:C2), keeping D2's note text intact.
$wsSingle = $wb.Worksheets.Item("Single-storm shps")
$wsSingle.Range("A2:C2").Clear()
$wsSingle.Range("A2:C2").Select()

# Make "Single-storm shps" the active sheet/tab, matching the saved view state.
$wsSingle.Activate()
